# Generate Report for Handoff
# Updates the localization-status report:
#  - Bumps the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#    for the rows that were just (re)handed off.
#  - Marks the Priority column as "ht" for those same rows on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 9, 10, 12, 13, 14)

foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-08-26 22:21:50"
    $zhcn.Range("H$r").Value = "2016-08-26 22:21:45"
    $dede.Range("H$r").Value = "2016-08-26 22:21:50"

    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"
}
